$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @{ Row = 13; A = 13; B = "Elisângela dos Santos"; C = 13; D = "Rua Apucarana, 11, Ouro Preto, Belo Horizonte, MG"; E = 0 },
    @{ Row = 14; A = 14; B = "Pedro Castro";          C = 14; D = "Rua Colorado 87";                                  E = 0 },
    @{ Row = 15; A = 15; B = "Clara Maria Paiva";     C = 15; D = "Rua Rio Grande 76";                                E = 0 },
    @{ Row = 16; A = 16; B = "Leila Martins";         C = 16; D = "Rua Itabirito 379";                                E = 0 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
